# Update cryptocurrency price/volume snapshot values (scraped refresh).
# Column D = Price (text), Column E = Volume(1h) % change (text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.636.63"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "'1.848.38"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").Value = "'1.033"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'321.68"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'1.028"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4377"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.3786"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "'0.07380"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.8808"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'21.48"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "'1.857.88"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "'5.493"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "'6.699"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "'0.07137"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'84.99"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "'1.033"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'0.000009043"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'1.027"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'15.43"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'27.667.22"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "'5.288"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "'2.087.87"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "'2.042"
$ws.Range("E25").Value = "  +6.30%  "
$ws.Range("D26").Value = "'157.67"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").Value = "'18.68"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'1.997"
$ws.Range("E28").Value = "  +3.18%  "
$ws.Range("D29").Value = "'5.327"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").Value = "'117.64"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").Value = "'0.09030"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "'0.7683"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("D34").Value = "'2.992"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("D35").Value = "'4.546"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").Value = "'1.028"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'1.143"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'0.01969"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'0.05257"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "'2.834"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'0.5174"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'0.1668"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'6.835"
$ws.Range("E43").Value = "  +3.29%  "
$ws.Range("D44").Value = "'8.748"
$ws.Range("E44").Value = "  +2.46%  "
$ws.Range("D45").Value = "'110.28"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'10.63"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'0.06615"
$ws.Range("E47").Value = "  +4.12%  "
$ws.Range("D48").Value = "'1.030"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'1.695"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "'0.4686"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "'1.888"
$ws.Range("E51").Value = "  -0.46%  "
